{"js": "// Word Add-in (Office.js) script \u2014 body of `async (context) => { ... }`.\n//\n// Applies three content edits to the \"Week \u2026\" addenda list, matching the\n// target OOXML from the diff:\n//   1. \"Computations (plus, minus, modulus, divide, multiplication,\u2026)\"\n//      (under Week 2) is re-run-split with spell/grammar w:proofErr marks\n//      around \"modulus\", \"divide\" and \"multiplication,\u2026\".\n//   2. \"Classes and Resources \u2013 remove, need to know operators for this\"\n//      (under Week 11) becomes simply \"Design\".\n//   3. \"\u2026 Records and Fields (CorporateInfo from Week 8)\" (under Week 12)\n//      is re-run-split with a w:proofErr spell mark around \"CorporateInfo\".\n//\n// Because w:proofErr markers are not exposed through the high-level\n// Office.js object model, each paragraph is rewritten in place with\n// Range.insertOoxml(..., Word.InsertLocation.replace) using a minimal\n// flat-OPC wrapper so the exact run / w:proofErr structure can be\n// produced.\n\nconst WORD_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction flatOpcDocument(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    \"<w:document \" + WORD_NS + \">\" +\n    \"<w:body>\" + bodyInnerXml + \"</w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nasync function replaceParagraphOoxml(context, paragraph, newParagraphInnerXml) {\n  const range = paragraph.getRange();\n  const xml = flatOpcDocument(\"<w:p>\" + newParagraphInnerXml + \"</w:p>\");\n  range.insertOoxml(xml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst computationsPara = paragraphs.items.find(\n  (p) => p.text === \"Computations (plus, minus, modulus, divide, multiplication,\\u2026)\"\n);\nconst classesPara = paragraphs.items.find(\n  (p) => p.text === \"Classes and Resources \\u2013 remove, need to know operators for this\"\n);\nconst secondaryStoragePara = paragraphs.items.find(\n  (p) => p.text === \"Secondary Storage (File I/O) \\u2013 Records and Fields (CorporateInfo from Week 8)\"\n);\n\nif (!computationsPara || !classesPara || !secondaryStoragePara) {\n  throw new Error(\"Could not locate one or more target paragraphs\");\n}\n\n// 1) \"Computations (plus, minus, modulus, divide, multiplication,\u2026)\"\nconst rPr = '<w:rPr><w:lang w:val=\"fr-FR\"/></w:rPr>';\nconst computationsInner =\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">Computations (plus, minus, </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + rPr + '<w:t>modulus</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">, </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + rPr + '<w:t>d</w:t></w:r>' +\n  '<w:r>' + rPr + '<w:t>ivide</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r>' + rPr + '<w:t xml:space=\"preserve\">, </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r>' + rPr + '<w:t>multiplication,\\u2026</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r>' + rPr + '<w:t>)</w:t></w:r>';\nawait replaceParagraphOoxml(context, computationsPara, computationsInner);\n\n// 2) \"Classes and Resources \u2013 remove, need to know operators for this\" -> \"Design\"\nconst classesInner = '<w:r><w:t>Design</w:t></w:r>';\nawait replaceParagraphOoxml(context, classesPara, classesInner);\n\n// 3) \"\u2026 Records and Fields (CorporateInfo from Week 8)\"\nconst secondaryStorageInner =\n  '<w:r><w:t>Secondary Storage (File I/O) \\u2013 Records and Fields</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>CorporateInfo</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> from Week 8)</w:t></w:r>';\nawait replaceParagraphOoxml(context, secondaryStoragePara, secondaryStorageInner);\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Applies three content edits to the \"Week \u2026\" addenda list, matching the\n# target OOXML from the diff:\n#   1. \"Computations (plus, minus, modulus, divide, multiplication,\u2026)\"\n#      (under Week 2) is re-run-split with spell/grammar w:proofErr marks\n#      around \"modulus\", \"divide\" and \"multiplication,\u2026\".\n#   2. \"Classes and Resources \u2013 remove, need to know operators for this\"\n#      (under Week 11) becomes simply \"Design\".\n#   3. \"\u2026 Records and Fields (CorporateInfo from Week 8)\" (under Week 12)\n#      is re-run-split with a w:proofErr spell mark around \"CorporateInfo\".\n#\n# w:proofErr markers aren't reachable through the Range/Paragraph text\n# properties, so each target paragraph is rewritten in place with\n# Range.InsertXML(...) using a minimal flat-OPC wrapper so the exact\n# run / w:proofErr structure from the diff can be produced.\n\n$d = $word.ActiveDocument\n\n$enDash = [char]0x2013\n$ellipsis = [char]0x2026\n\nfunction Get-ParagraphByText($doc, [string]$text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    throw \"Paragraph not found: $text\"\n}\n\nfunction New-FlatOpcXml([string]$paragraphInnerXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' + $paragraphInnerXml + '</w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\n# 1) \"Computations (plus, minus, modulus, divide, multiplication,\u2026)\"\n$computationsText = \"Computations (plus, minus, modulus, divide, multiplication,$ellipsis)\"\n$computationsPara = Get-ParagraphByText $d $computationsText\n\n$rPr = '<w:rPr><w:lang w:val=\"fr-FR\"/></w:rPr>'\n$computationsInner =\n    \"<w:r>$rPr<w:t xml:space=`\"preserve`\">Computations (plus, minus, </w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r>$rPr<w:t>modulus</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"<w:r>$rPr<w:t xml:space=`\"preserve`\">, </w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r>$rPr<w:t>d</w:t></w:r>\" +\n    \"<w:r>$rPr<w:t>ivide</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"<w:r>$rPr<w:t xml:space=`\"preserve`\">, </w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r>$rPr<w:t>multiplication,&#8230;</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"<w:r>$rPr<w:t>)</w:t></w:r>\"\n[void]$computationsPara.Range.InsertXML((New-FlatOpcXml $computationsInner))\n\n# 2) \"Classes and Resources \u2013 remove, need to know operators for this\" -> \"Design\"\n$classesText = \"Classes and Resources $enDash remove, need to know operators for this\"\n$classesPara = Get-ParagraphByText $d $classesText\n$classesInner = '<w:r><w:t>Design</w:t></w:r>'\n[void]$classesPara.Range.InsertXML((New-FlatOpcXml $classesInner))\n\n# 3) \"\u2026 Records and Fields (CorporateInfo from Week 8)\"\n$secondaryStorageText = \"Secondary Storage (File I/O) $enDash Records and Fields (CorporateInfo from Week 8)\"\n$secondaryStoragePara = Get-ParagraphByText $d $secondaryStorageText\n$secondaryStorageInner =\n    \"<w:r><w:t>Secondary Storage (File I/O) &#8211; Records and Fields</w:t></w:r>\" +\n    \"<w:r><w:t xml:space=`\"preserve`\"> (</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>CorporateInfo</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"<w:r><w:t xml:space=`\"preserve`\"> from Week 8)</w:t></w:r>\"\n[void]$secondaryStoragePara.Range.InsertXML((New-FlatOpcXml $secondaryStorageInner))\n"}
